$d = $word.ActiveDocument

# ===========================================================================
# NOTE on paragraph indices below: they refer to the *original* document's
# paragraph numbering (1-based).  Edits are applied from the end of the
# document towards the beginning so that earlier, not-yet-processed
# paragraph indices stay valid while later ones are being restructured.
# ===========================================================================

# ---------------------------------------------------------------------------
# 1) "Next Steps" section (originally paragraphs 43-51)
# ---------------------------------------------------------------------------

# 1a) paragraph 43: rewrite sentence
$oldA = "In order to validate the usefulness of this idea, we should partner with FEMA or a local government to use VAN in a disaster training exercise."
$newA = "As a first step to validating this idea, we can partner with FEMA or a local government to use VAN in a disaster response training or simulation. Medium and large cities perform these types of simulations regularly. "
$d.Content.Find.Execute($oldA, $true, $false, $false, $false, $false, $true, 1, $false, $newA, 2) | Out-Null

# 1b) paragraph 44: change opening words
$oldB = "Get approval and"
$newB = "We should seek approval to"
$d.Content.Find.Execute($oldB, $true, $false, $false, $false, $false, $true, 1, $false, $newB, 2) | Out-Null

# 1c) remove paragraphs 45 (blank), 46 ("-Proof of concept..."), 47
#     ("-Pros/Cons...FEMA"), 48 ("-Pros/Cons...local government"), and 49
#     (first of the three trailing blank/bold paragraphs); then insert a new
#     paragraph with the "To test the efficacy..." text where they were.
$p44 = $d.Paragraphs.Item(44)
$p49 = $d.Paragraphs.Item(49)
$delRange = $d.Range($p44.Range.End, $p49.Range.End)
$delRange.Delete()

$p44b = $d.Paragraphs.Item(44)
$p44b.Range.InsertParagraphAfter()
$p45new = $d.Paragraphs.Item(45)
$p45new.Range.Text = "To test the efficacy of the product, we can work with the MIT Humanitarian Response Lab to design a study assessing the impact of VAN software during the simulation compared to other alternatives. "
$p45new.Format.FirstLineIndent = 36

# ---------------------------------------------------------------------------
# 2) "Customizations for Disaster Recovery" section (originally paragraphs
#    37-40): merge the heading, "For a serious roll out..." paragraph, a
#    blank paragraph, and "Existing processes and in house solutions" into
#    one paragraph.
# ---------------------------------------------------------------------------
$p37 = $d.Paragraphs.Item(37)
$insPos = $p37.Range.End - 1
$d.Range($insPos, $insPos).InsertBefore(": ")
$p37b = $d.Paragraphs.Item(37)
$spacePos = $p37b.Range.End - 2
$d.Range($spacePos, $spacePos + 1).Underline = 0

$oldRollout = "For a serious roll out, customizations may at least include "
$newRollout = "For a serious roll out, customizations may at least need to include language changes to reflect disaster recovery and not social organizing. "
$d.Content.Find.Execute($oldRollout, $true, $false, $false, $false, $false, $true, 1, $false, $newRollout, 2) | Out-Null

$p40 = $d.Paragraphs.Item(40)
$existingTextRange = $d.Range($p40.Range.Start, $p40.Range.End - 1)
$existingTextRange.Delete()

$p37c = $d.Paragraphs.Item(37)
$d.Range($p37c.Range.End - 1, $p37c.Range.End).Delete()
$p37d = $d.Paragraphs.Item(37)
$d.Range($p37d.Range.End - 1, $p37d.Range.End).Delete()
$p37e = $d.Paragraphs.Item(37)
$d.Range($p37e.Range.End - 1, $p37e.Range.End).Delete()

# ---------------------------------------------------------------------------
# 3) Remove one of the two blank paragraphs before the "Customizations"
#    heading (originally paragraphs 35 and 36; paragraph 35 is merged away).
# ---------------------------------------------------------------------------
$p35 = $d.Paragraphs.Item(35)
$d.Range($p35.Range.End - 1, $p35.Range.End).Delete()

# ---------------------------------------------------------------------------
# 4) "Data and Administrative Ownership" section (originally paragraphs
#    30-31): add a trailing colon and merge into the "While VAN makes..."
#    paragraph.
# ---------------------------------------------------------------------------
$oldOwn = "Data and Administrative Ownership "
$newOwn = "Data and Administrative Ownership:"
$d.Content.Find.Execute($oldOwn, $true, $false, $false, $false, $false, $true, 1, $false, $newOwn, 2) | Out-Null

$p30 = $d.Paragraphs.Item(30)
$ownInsPos = $p30.Range.End - 1
$d.Range($ownInsPos, $ownInsPos).InsertBefore(" ")
$p30b = $d.Paragraphs.Item(30)
$ownSpacePos = $p30b.Range.End - 2
$d.Range($ownSpacePos, $ownSpacePos + 1).Underline = 0

$p30c = $d.Paragraphs.Item(30)
$d.Range($p30c.Range.End - 1, $p30c.Range.End).Delete()

# ---------------------------------------------------------------------------
# 5) "Training and existing trained user base:" section (originally
#    paragraphs 25-26): insert an underlined trailing space and merge into
#    the following paragraph.
# ---------------------------------------------------------------------------
$p25 = $d.Paragraphs.Item(25)
$trainInsPos = $p25.Range.End - 1
$d.Range($trainInsPos, $trainInsPos).InsertBefore(" ")

$p25b = $d.Paragraphs.Item(25)
$d.Range($p25b.Range.End - 1, $p25b.Range.End).Delete()
